$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases & Results")

# --- Row 45 (TestCase_ID 43, REQ-42) ---
$ws.Range("E45").Value = "Mid Impact"
$ws.Range("F45").Value = "Test that The app waits for the response and display the description"
$ws.Range("G45").Value = "User must pick a book from woodlands library"
$ws.Range("H45").Value = "Pick any book from the list of available book In woodlands and press its item then select AI"
$ws.Range("I45").Value = "A description of the Book should be shown on the page"
$ws.Range("J45").Value = "A description of the Book shown on the page"
$ws.Rows.Item(45).RowHeight = 57.6

# --- Row 46 (TestCase_ID 44, REQ-43) ---
$ws.Range("E46").Value = "High Impact"
$ws.Range("F46").Value = "Test that Once user confirms, App should check if the user has borrowed 10 books already"
$ws.Range("G46").Value = "User must pick a book from woodlands library"
$ws.Range("H46").Value = "Pick any book from the list of available book In woodlands and press its item then select Borrow"
$ws.Range("I46").Value = "No visible output, go to next test case"
$ws.Range("J46").Value = "No visible output, go to next test case"
$ws.Rows.Item(46).RowHeight = 57.6

# --- Row 47 (TestCase_ID 45, REQ-44) ---
$ws.Range("E47").Value = "High Impact"
$ws.Range("F47").Value = "If the user has borrowed more than 10 books from REQ-43, the app should reject the reservation"
$ws.Range("G47").Value = "User must pick a book from woodlands library and borrow it"
$ws.Range("H47").Value = "Borrow 10 books through the previous steps then try to borrow an 11th one"
$ws.Range("I47").Value = "Snackbar appears and disallow the reservation"
$ws.Range("J47").Value = "Snackbar appears and disallow the reservation"
$ws.Rows.Item(47).RowHeight = 43.2

# --- Update the view/selection to match the edited area ---
$ws.Activate()
$ws.Range("F46").Select()
